{"js": "// The target run (bold, noProof) currently reads:\n//   \"only for bounding box that is responsible for detecting the object.\"\n// It needs to become:\n//   \"only for the bounding box that is responsible for detecting the object.\"\n// i.e. insert the word \"the \" right after \"only for \".\nconst body = context.document.body;\n\nconst results = body.search(\"only for bounding box that is responsible for detecting the object.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to edit.\");\n}\n\n// There should only be a single match in the document; search within it for the\n// exact insertion anchor so the edit lands precisely between \"only for \" and\n// \"bounding box...\".\nconst target = results.items[0];\nconst subResults = target.search(\"only for \", { matchCase: true });\nsubResults.load(\"items\");\nawait context.sync();\n\nif (subResults.items.length === 0) {\n  throw new Error(\"Could not find 'only for ' inside the target sentence.\");\n}\n\n// Collapse to a caret right after \"only for \" and insert \"the \" there,\n// preserving the run's existing bold/noProof formatting.\nconst afterOnlyFor = subResults.items[0].getRange(\"End\");\nafterOnlyFor.insertText(\"the \", \"Before\");\nawait context.sync();\n", "ps1": "# The target run (bold, noProof) currently reads:\n#   \"only for bounding box that is responsible for detecting the object.\"\n# It needs to become:\n#   \"only for the bounding box that is responsible for detecting the object.\"\n# i.e. insert the word \"the \" right after \"only for \", keeping the existing\n# bold / noProof character formatting intact.\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"only for bounding box that is responsible for detecting the object.\"\n$find.Execute() | Out-Null\n\nif ($find.Found) {\n    # Narrow to just \"only for \" within the matched sentence so the insert\n    # lands precisely before \"bounding box\", without retyping the rest of\n    # the (already correctly formatted) sentence.\n    $insertAt = $range.Duplicate\n    $insertAt.Collapse(1) | Out-Null  # wdCollapseStart\n    $insertAt.MoveEnd(1, 9) | Out-Null  # wdCharacter, len(\"only for \") == 9\n    $insertAt.Collapse(0) | Out-Null  # wdCollapseEnd\n    $insertAt.InsertBefore(\"the \") | Out-Null\n}\n"}
